$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new (older) observation was added to the series, so a row is inserted
# right after the header, shifting all existing data rows down by one.
$ws.Rows.Item(2).Insert()

# The inserted row picks up blended formatting from the header row;
# reset it back to "no explicit style" like the other data rows, then
# restore the date style (s="2") on column A only, matching the rest
# of that column.
$ws.Range("B2:E2").ClearFormats()
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New row 2 values (oldest data point, now present in the series).
$ws.Range("A2").Value = 39400
$ws.Range("B2").Value = 2007
$ws.Range("D2").Value = 2008

# The forecast columns (C: y_0_forecast, E: y_1_forecast) were recomputed
# for every data row (now rows 2..19) following the bugfix.
$values = @(
    @(2, 4.930115226412335, 0.3630458632513767),
    @(3, 1.457587285166628, 0.507956838644974),
    @(4, -0.9140166223623569, 5.6395352704941),
    @(5, 2.585942866987878, 1.724360951547554),
    @(6, 4.253963781362402, -0.6955733540840336),
    @(7, 1.752870900283909, 4.300339264728548),
    @(8, -1.479696720105139, 8.296896928314457),
    @(9, 3.900127535411246, -2.092856741436233),
    @(10, 0.03947433952959933, -1.259568900987018),
    @(11, 2.192778679161944, 0.5033587260849126),
    @(12, 3.40836448860673, -0.3010260522302355),
    @(13, 2.799070570134488, 3.825329033908775),
    @(14, 4.195393191694419, 3.942709467505678),
    @(15, 1.666553973046048, -0.2638638106667313),
    @(16, 1.879266440112803, -0.04532879466145889),
    @(17, -2.620683231370946, -3.179374983142691),
    @(18, -3.036556262700274, -2.51939929628594),
    @(19, -2.953443685011514, -2.112604539331953)
)

foreach ($row in $values) {
    $r = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 5).Value = $row[2]
}
